# Daily attendance processing - 2025-11-04 02:59:29
# Reverses the order of the comma-separated "Recorded By" entries in column G
# for every data row on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Text

    if ([string]::IsNullOrEmpty($val)) {
        continue
    }

    $parts = $val -split ',\s*'

    if ($parts.Count -le 1) {
        continue
    }

    $reversed = $parts[($parts.Count - 1)..0]
    $newVal = $reversed -join ', '

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
